# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.420.50"
$ws.Range("E2").Value = "  -0.91%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.773.93"
$ws.Range("E3").Value = "  -0.22%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "353.57"
$ws.Range("E5").Value = "  -0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.84"
$ws.Range("E6").Value = "  -1.75%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.549"
$ws.Range("E7").Value = "  -2.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.585"
$ws.Range("E9").Value = "  -1.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.66"
$ws.Range("E10").Value = "  -1.22%  "

$ws.Range("E11").Value = "  +3.18%  "

$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0832"
$ws.Range("E12").Value = "  -2.37%  "

$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.00"
$ws.Range("E13").Value = "  +3.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.55"
$ws.Range("E14").Value = "  -0.92%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.210.26"
$ws.Range("E15").Value = "  -0.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.774.84"
$ws.Range("E16").Value = "  +0.26%  "

$ws.Range("E17").Value = "  -0.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.379.83"
$ws.Range("E18").Value = "  -0.76%  "

$ws.Range("E19").Value = "  +2.40%  "

$ws.Range("E20").Value = "  -2.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.14"
$ws.Range("E21").Value = "  +0.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0961"
$ws.Range("E22").Value = "  -1.63%  "

$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "265.47"
$ws.Range("E24").Value = "  -3.17%  "

$ws.Range("E25").Value = "  -1.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.98"
$ws.Range("E27").Value = "  -2.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.160"
$ws.Range("E28").Value = "  +11.98%  "

$ws.Range("E29").Value = "  +0.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.66"
$ws.Range("E30").Value = "  +7.82%  "

$ws.Range("E31").Value = "  -0.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.18"
$ws.Range("E32").Value = "  +8.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "51.88"
$ws.Range("E33").Value = "  +0.76%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0447"
$ws.Range("E34").Value = "  -3.21%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.50"
$ws.Range("E35").Value = "  +5.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0825"
$ws.Range("E36").Value = "  -2.41%  "

$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.21"
$ws.Range("E38").Value = "  +1.02%  "

$ws.Range("E39").Value = "  -2.46%  "

$ws.Range("E40").Value = "  -1.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.54"
$ws.Range("E41").Value = "  +0.70%  "

$ws.Range("E42").Value = "  -0.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "120.90"
$ws.Range("E43").Value = "  -0.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.97"
$ws.Range("E44").Value = "  -0.21%  "

$ws.Range("E45").Value = "  -1.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.094.28"
$ws.Range("E46").Value = "  +1.37%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.24"
$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.32"
$ws.Range("E48").Value = "  +3.75%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.41"
$ws.Range("E49").Value = "  -4.73%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.900"
$ws.Range("E50").Value = "  -2.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.31"
$ws.Range("E51").Value = "  +8.07%  "
